$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.065.63"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "2.563.26"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'584.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.88%  "

$ws.Range("D6").Value = "'147.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +1.60%  "

$ws.Range("E9").Value = "  +2.88%  "

$ws.Range("E10").Value = "  +0.84%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("D13").Value = "'27.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.53%  "

$ws.Range("D14").Value = "3.022.11"
$ws.Range("E14").Value = "  +0.72%  "

$ws.Range("D15").Value = "63.025.02"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("E16").Value = "  +2.12%  "

$ws.Range("D17").Value = "2.567.56"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").Value = "'344.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.42%  "

$ws.Range("E20").Value = "  +2.78%  "

$ws.Range("E21").Value = "  +1.61%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("E23").Value = "  -3.70%  "

$ws.Range("D24").Value = "'66.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.99%  "

$ws.Range("D25").Value = "2.693.63"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("D28").Value = "'8.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.45%  "

$ws.Range("E29").Value = "  +0.39%  "

$ws.Range("E30").Value = "  -1.71%  "

$ws.Range("D31").Value = "'8.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.83%  "

$ws.Range("D32").Value = "'1.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.05%  "

$ws.Range("D33").Value = "0.0₃0822"
$ws.Range("E33").Value = "  -0.18%  "

$ws.Range("D34").Value = "'460.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.97%  "

$ws.Range("D35").Value = "'175.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "

$ws.Range("E36").Value = "  +2.56%  "

$ws.Range("D37").Value = "'0.406"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.00%  "

$ws.Range("D38").Value = "'19.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.88%  "

$ws.Range("E39").Value = "  +3.06%  "

$ws.Range("E41").Value = "  -0.71%  "

$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").Value = "'150.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("E44").Value = "  +1.19%  "

$ws.Range("E45").Value = "  -0.51%  "

$ws.Range("D46").Value = "'0.0545"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.60%  "

$ws.Range("E47").Value = "  +1.21%  "

$ws.Range("D48").Value = "'0.0973"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.05%  "

$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("E50").Value = "  -2.57%  "

$ws.Range("D51").Value = "'11.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.61%  "
